$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Title cell A1 -> placeholder "{{Name}}" (banner style/fill/merge untouched here).
$ws.Range("A1").Value = "{{Name}}"

# 2. New column J absorbed into the report: give each row-1/2/3 cell the same
#    look as its row neighbour (format-only paste, so no stray style gets
#    minted), then write its own text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J2").Value = "Lý do từ chối"

$ws.Range("E3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value = "{{Items.Notes}}"

# 3. F3 swaps its old vertical-center-only look for the plain centered
#    data-cell look shared by the rest of row 3 (value is unchanged).
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "{{Items.Year}}"

# 4. Re-merge the title bar across the new column.
$ws.Range("A1:I1").UnMerge()
$ws.Range("A1:J1").Merge()

# 5. Give column J the same kind of width as the other wide text column.
$ws.Columns.Item(10).ColumnWidth = 47

# 6. Selection as last left by the author.
$ws.Range("F3").Select()
